$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Romania Liga I")

# ---------------------------------------------------------------------------
# Rows 235-238: the four fixtures (everything except columns A, C, D, E, I)
# get reordered - row 235's data moves to row 238, 236 <-> 237, and 238's
# data moves to row 235. Columns A (index), C/D (div names), E (date) and I
# (FTAG) stay put on their row. Apply the new values directly.
# ---------------------------------------------------------------------------

# New row 235 (was row 238's data)
$ws.Range("B235").Value = 6861095
$ws.Range("F235").Value = "FC Botosani"
$ws.Range("G235").Value = "Farul Constanta"
$ws.Range("H235").Value = 0
$ws.Range("J235").Value = "D"
$ws.Range("K235").Value = 3.75
$ws.Range("L235").Value = 3.4
$ws.Range("M235").Value = 1.909
$ws.Range("N235").Value = 3.1
$ws.Range("O235").Value = 3
$ws.Range("P235").Value = 2.375
$ws.Range("Q235").Value = 0.25
$ws.Range("R235").Value = 1.775
$ws.Range("S235").Value = 2.1
$ws.Range("T235").Value = 2
$ws.Range("U235").Value = 1.8
$ws.Range("V235").Value = 2.05
$ws.Range("W235").Value = -1
$ws.Range("X235").Value = 2
$ws.Range("Y235").Value = -1
$ws.Range("Z235").Value = 0.3875
$ws.Range("AA235").Value = -0.5
$ws.Range("AB235").Value = -1
$ws.Range("AC235").Value = 1.05

# New row 236 (was row 237's data)
$ws.Range("B236").Value = 6865915
$ws.Range("F236").Value = "FC Voluntari"
$ws.Range("G236").Value = "Universitatea Cluj"
$ws.Range("H236").Value = 0
$ws.Range("J236").Value = "D"
$ws.Range("K236").Value = 3.5
$ws.Range("L236").Value = 3.25
$ws.Range("M236").Value = 2.05
$ws.Range("N236").Value = 3.4
$ws.Range("O236").Value = 3.1
$ws.Range("P236").Value = 2.15
$ws.Range("Q236").Value = 0.25
$ws.Range("R236").Value = 1.975
$ws.Range("S236").Value = 1.875
$ws.Range("T236").Value = 2.25
$ws.Range("U236").Value = 2.05
$ws.Range("V236").Value = 1.75
$ws.Range("W236").Value = -1
$ws.Range("X236").Value = 2.1
$ws.Range("Y236").Value = -1
$ws.Range("Z236").Value = 0.4875
$ws.Range("AA236").Value = -0.5
$ws.Range("AB236").Value = -1
$ws.Range("AC236").Value = 0.75

# New row 237 (was row 236's data)
$ws.Range("B237").Value = 6836277
$ws.Range("F237").Value = "CFR Cluj"
$ws.Range("G237").Value = "AFC Hermannstadt"
$ws.Range("H237").Value = 1
$ws.Range("J237").Value = "H"
$ws.Range("K237").Value = 1.7
$ws.Range("L237").Value = 3.4
$ws.Range("M237").Value = 5
$ws.Range("N237").Value = 1.65
$ws.Range("O237").Value = 3.5
$ws.Range("P237").Value = 5.25
$ws.Range("Q237").Value = -0.75
$ws.Range("R237").Value = 1.85
$ws.Range("S237").Value = 2
$ws.Range("T237").Value = 2.25
$ws.Range("U237").Value = 1.875
$ws.Range("V237").Value = 1.975
$ws.Range("W237").Value = 0.6499999999999999
$ws.Range("X237").Value = -1
$ws.Range("Y237").Value = -1
$ws.Range("Z237").Value = 0.425
$ws.Range("AA237").Value = -0.5
$ws.Range("AB237").Value = -1
$ws.Range("AC237").Value = 0.9750000000000001

# New row 238 (was row 235's data)
$ws.Range("B238").Value = 6852370
$ws.Range("F238").Value = "Dinamo Bucharest"
$ws.Range("G238").Value = "ACS UTA Batrana Doamna"
$ws.Range("H238").Value = 1
$ws.Range("J238").Value = "H"
$ws.Range("K238").Value = 2.55
$ws.Range("L238").Value = 2.875
$ws.Range("M238").Value = 3
$ws.Range("N238").Value = 2.375
$ws.Range("O238").Value = 3
$ws.Range("P238").Value = 3.1
$ws.Range("Q238").Value = -0.25
$ws.Range("R238").Value = 2
$ws.Range("S238").Value = 1.85
$ws.Range("T238").Value = 2.25
$ws.Range("U238").Value = 1.975
$ws.Range("V238").Value = 1.875
$ws.Range("W238").Value = 1.375
$ws.Range("X238").Value = -1
$ws.Range("Y238").Value = -1
$ws.Range("Z238").Value = 1
$ws.Range("AA238").Value = -1
$ws.Range("AB238").Value = -1
$ws.Range("AC238").Value = 0.875

# ---------------------------------------------------------------------------
# Rows 259-264: odds updates for upcoming fixtures (simple value refreshes).
# ---------------------------------------------------------------------------

# Row 259
$ws.Range("R259").Value = 2.025
$ws.Range("S259").Value = 1.825
$ws.Range("U259").Value = 1.9
$ws.Range("V259").Value = 1.95

# Row 260
$ws.Range("U260").Value = 1.825
$ws.Range("V260").Value = 2.025

# Row 261
$ws.Range("N261").Value = 2.25
$ws.Range("O261").Value = 3.1
$ws.Range("P261").Value = 3.1
$ws.Range("R261").Value = 2
$ws.Range("S261").Value = 1.85

# Row 262
$ws.Range("N262").Value = 2.05
$ws.Range("O262").Value = 3.3
$ws.Range("P262").Value = 3.75
$ws.Range("Q262").Value = -0.5
$ws.Range("R262").Value = 2.05
$ws.Range("S262").Value = 1.8

# Row 264
$ws.Range("N264").Value = 1.7
$ws.Range("O264").Value = 3.6
$ws.Range("P264").Value = 4.75
$ws.Range("Q264").Value = -0.75
$ws.Range("R264").Value = 2
$ws.Range("S264").Value = 1.85
$ws.Range("U264").Value = 1.875
$ws.Range("V264").Value = 1.975
